# add unit in bulk
# Remove the "SKU" column (column A) from the add-unit template, shifting
# id_wh / serial_number / comment left by one column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely (this removes the SKU / "Item's SKU" header+value
# and shifts B:D -> A:C).
$ws.Columns.Item(1).Delete()

# Select column B (matches the resulting selection in the saved file).
$ws.Range("B1:B1048576").Select()
